# edit.ps1 - applies the diff described in the task to the document.
#
# Strategy:
#  - For simple "merge adjacent runs into one run with identical text"
#    edits (headings, table-cell merges), use Find/Replace: Word (and this
#    runtime) collapses the matched range - even if it spans several runs -
#    into a single run using the first run's formatting.
#  - For edits that split a single run's text into several runs with
#    identical run-properties (rPr), Find/Replace can't produce that
#    shape (it always collapses to one run), so we use Range.InsertXML
#    with a hand-built WordOpenXML package fragment that lists the exact
#    <w:r> runs (and occasionally <w:proofErr/> siblings) we need. The
#    paragraph's own identity attributes (w14:paraId/w14:textId/w:rsidR/...)
#    and <w:pPr> are always re-supplied explicitly in the fragment so nothing
#    is lost (table-cell paragraphs do not auto-preserve them the way body
#    paragraphs do).
#  - Paragraph deletions use Paragraph.Range.Delete(), which removes the
#    paragraph mark too and merges cleanly with the following paragraph.

$d = $word.ActiveDocument

function Get-ParaIndexByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text
        # Trim table "end of cell" markers (CR + BEL) and plain trailing CR
        # that Range.Text reports but that aren't part of the visible text.
        $t2 = $t.TrimEnd([char]13, [char]7)
        if ($t2 -eq $text) {
            return $i
        }
    }
    return -1
}

function Set-ParagraphRuns($doc, $paraText, $paraXmlInner) {
    # Replaces the *content* (everything except the paragraph mark) of the
    # paragraph whose flattened text equals $paraText with the raw OOXML
    # given in $paraXmlInner (one or more sibling nodes: <w:pPr>, <w:r>,
    # <w:proofErr/>, ...). The supplied fragment must explicitly restate the
    # original paragraph's w14:paraId/w14:textId/w:rsidR/w:rsidRDefault (and
    # w:rsidP if present) since they are not guaranteed to be preserved
    # automatically inside table cells.
    $idx = Get-ParaIndexByText $doc $paraText
    if ($idx -lt 0) {
        throw "Paragraph not found: $paraText"
    }
    $p = $doc.Paragraphs.Item($idx)
    $full = $p.Range
    $s = $full.Start
    $e = $full.End
    $inner = $doc.Range($s, $e - 1)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + $paraXmlInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $inner.InsertXML($xml)
}

function Replace-Text($doc, $findText, $replaceText) {
    $found = $doc.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $found) {
        throw "Find text not found: $findText"
    }
}

# ---------------------------------------------------------------------
# 1) "No explicit go-to-market or beachhead customer" -> four runs that
#    read "Financial forecast is overly used in traditional corporate
#    setting, but I think it doesn't mean much in venture building. "
# ---------------------------------------------------------------------
$rightQuote = [char]0x2019
$run3Text = "s overly used in traditional corporate setting, but I think it doesn" + $rightQuote + "t mean much in venture building."

$inner1 = '<w:p w14:paraId="4760650C" w14:textId="77777777" w:rsidR="003002E6" w:rsidRDefault="00000000">' + `
    '<w:pPr><w:pStyle w:val="ListBullet"/></w:pPr>' + `
    '<w:r><w:rPr><w:rFonts w:ascii="IBM Plex Sans" w:eastAsia="IBM Plex Sans" w:hAnsi="IBM Plex Sans"/></w:rPr><w:t>F</w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:ascii="IBM Plex Sans" w:eastAsia="IBM Plex Sans" w:hAnsi="IBM Plex Sans"/></w:rPr><w:t>inancial forecast i</w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:ascii="IBM Plex Sans" w:eastAsia="IBM Plex Sans" w:hAnsi="IBM Plex Sans"/></w:rPr><w:t>' + $run3Text + '</w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:ascii="IBM Plex Sans" w:eastAsia="IBM Plex Sans" w:hAnsi="IBM Plex Sans"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
    '</w:p>'
Set-ParagraphRuns $d "No explicit go-to-market or beachhead customer" $inner1

# ---------------------------------------------------------------------
# 2) "Partner governance and IP protection unclear" -> single run:
#    "Revenue is too dependent on partner cost and is tied to amount of
#    work done, not value created. "
# ---------------------------------------------------------------------
Replace-Text $d "Partner governance and IP protection unclear" "Revenue is too dependent on partner cost and is tied to amount of work done, not value created. "

# ---------------------------------------------------------------------
# 3) Remove the two bullets "10x revenue growth assumptions aggressive"
#    and "'Execution intelligence' not operationally defined", plus the
#    blank paragraph right after them.
# ---------------------------------------------------------------------
$idx = Get-ParaIndexByText $d "10x revenue growth assumptions aggressive"
if ($idx -lt 0) { throw "paragraph not found: 10x revenue growth..." }
$d.Paragraphs.Item($idx).Range.Delete()

$idx = Get-ParaIndexByText $d "'Execution intelligence' not operationally defined"
if ($idx -lt 0) { throw "paragraph not found: 'Execution intelligence'..." }
$d.Paragraphs.Item($idx).Range.Delete()

$idx = Get-ParaIndexByText $d ""
# The blank paragraph right after is the one immediately preceding
# "3. Side-by-Side Comparison: AVEP vs OSR" - locate it precisely instead
# of relying on the first blank paragraph in the whole document.
$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "") {
        $nextT = $d.Paragraphs.Item($i + 1).Range.Text.TrimEnd([char]13, [char]7)
        if ($nextT -eq "3. Side-by-Side Comparison: AVEP vs OSR") {
            $targetIdx = $i
            break
        }
    }
}
if ($targetIdx -lt 0) { throw "blank paragraph before section 3 not found" }
$d.Paragraphs.Item($targetIdx).Range.Delete()

# ---------------------------------------------------------------------
# 4) "Research automation, hypothesis generation" -> two runs:
#    "Research automation, " + "rapid prototyping and accelerated learnings "
# ---------------------------------------------------------------------
$inner4 = '<w:p w14:paraId="1E501C86" w14:textId="77777777" w:rsidR="003002E6" w:rsidRDefault="00000000">' + `
    '<w:r><w:rPr><w:rFonts w:ascii="IBM Plex Sans" w:eastAsia="IBM Plex Sans" w:hAnsi="IBM Plex Sans"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">Research automation, </w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:ascii="IBM Plex Sans" w:eastAsia="IBM Plex Sans" w:hAnsi="IBM Plex Sans"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">rapid prototyping and accelerated learnings </w:t></w:r>' + `
    '</w:p>'
Set-ParagraphRuns $d "Research automation, hypothesis generation" $inner4

# ---------------------------------------------------------------------
# 5) "Founders, solo entrepreneurs, startup teams" -> three runs plus
#    proofErr markers:
#    "Founders, entrepreneurs, startup teams" + ", angel investors, VC " +
#    (gramStart) "firms,.." (gramEnd)
# ---------------------------------------------------------------------
$inner5 = '<w:p w14:paraId="6EA340BF" w14:textId="77777777" w:rsidR="003002E6" w:rsidRDefault="00000000">' + `
    '<w:r><w:rPr><w:rFonts w:ascii="IBM Plex Sans" w:eastAsia="IBM Plex Sans" w:hAnsi="IBM Plex Sans"/><w:sz w:val="20"/></w:rPr><w:t>Founders, entrepreneurs, startup teams</w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:ascii="IBM Plex Sans" w:eastAsia="IBM Plex Sans" w:hAnsi="IBM Plex Sans"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">, angel investors, VC </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:rPr><w:rFonts w:ascii="IBM Plex Sans" w:eastAsia="IBM Plex Sans" w:hAnsi="IBM Plex Sans"/><w:sz w:val="20"/></w:rPr><w:t>firms,..</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '</w:p>'
Set-ParagraphRuns $d "Founders, solo entrepreneurs, startup teams" $inner5

# ---------------------------------------------------------------------
# 6) "Accumulated execution intelligence" + " + Storytelling" (two runs,
#    only in the OSR column) -> single run
#    "Accumulated execution intelligence + Storytelling"
# ---------------------------------------------------------------------
Replace-Text $d "Accumulated execution intelligence + Storytelling" "Accumulated execution intelligence + Storytelling"

# ---------------------------------------------------------------------
# 7) Heading run-merges: "N" + ".<rest>" -> single run "N.<rest>"
# ---------------------------------------------------------------------
Replace-Text $d "4. Critical Assessment: Which Model When?" "4. Critical Assessment: Which Model When?"
Replace-Text $d "4.1 AVEP is Better When:" "4.1 AVEP is Better When:"
Replace-Text $d "4.2 OSR is Better When:" "4.2 OSR is Better When:"
Replace-Text $d "5. Potential Synthesis: Hybrid Model" "5. Potential Synthesis: Hybrid Model"
Replace-Text $d "6. Fundamental Difference Summary" "6. Fundamental Difference Summary"
Replace-Text $d "8. Conclusion" "8. Conclusion"

# ---------------------------------------------------------------------
# 8) "Moderate" + " (cheap experiments" + " + seed investments" + ")"
#    -> single run "Moderate (cheap experiments + seed investments)"
# ---------------------------------------------------------------------
Replace-Text $d "Moderate (cheap experiments + seed investments)" "Moderate (cheap experiments + seed investments)"

# ---------------------------------------------------------------------
# 9) "7" + ". Recommendations" + " -" + " For Both" -> single run
#    "7. Recommendations"
# ---------------------------------------------------------------------
Replace-Text $d "7. Recommendations - For Both" "7. Recommendations"

Write-Output "Done"
